$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the Korean sample names with English placeholders (locale update)
$ws.Range("A2").Value = "Don Jone"
$ws.Range("A3").Value = "Jane"

# Move the frozen-pane (right side) selection to E10
$ws.Range("E10").Select() | Out-Null
